$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text / string cells -----------------------------------------------
# J2 ("004" -> "001") looks numeric, so Excel would otherwise coerce it to
# the number 1. Force a Text number format while writing it, then clear
# the format again so the cell keeps its original (default/General) style.
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "001"
$ws.Range("J2").ClearFormats()

# M2 / N2 contain date-like strings with a time component, which are kept
# as plain text automatically.
$ws.Range("M2").Value = "2020-12-24 00:00:00"
$ws.Range("N2").Value = "2019-12-31 00:00:00"

# --- Numeric cells -------------------------------------------------------
$ws.Range("O2").Value = 57553819.08
$ws.Range("P2").Value = 446235472.01
$ws.Range("Q2").Value = 389999005.12
$ws.Range("R2").Value = 19.1975545121
$ws.Range("S2").Value = 271347423.17
$ws.Range("T2").Value = 271347423.17
$ws.Range("U2").Value = 25.7791097276
$ws.Range("V2").Value = 41522510.4
$ws.Range("W2").Value = 27948094.87
$ws.Range("X2").Value = 2322600.43
$ws.Range("Y2").Value = 62675558.63
$ws.Range("Z2").Value = 63629216.28
$ws.Range("AA2").Value = 6075397.2
$ws.Range("AG2").Value = 3903528.56
$ws.Range("AP2").Value = 25.0106717829
$ws.Range("AQ2").Value = 49.235693191631
$ws.Range("AR2").Value = 63.130658129935
$ws.Range("AS2").Value = 59431719.08
$ws.Range("AT2").Value = 465.456482522652
